$d = $word.ActiveDocument

function Replace-TextPreservingRuns($OldText, $NewText, $RunPr, $SearchFrom) {
    $rng = $d.Content
    $rng.Start = $SearchFrom
    $found = $rng.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output ("NOT FOUND: " + $OldText)
        return $false
    }

    $target = $d.Range($rng.Start, $rng.End)

    $escaped = $NewText.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")

    $runXml = "<w:r>" + $RunPr + "<w:t>" + $escaped + "</w:t></w:r>"
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $runXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $target.InsertXML($xml)
    return $true
}

# 1. Title heading (first occurrence, no preceding empty run)
Replace-TextPreservingRuns "Play Juicy Reels for Free - Exciting Bonus Features" "Play Juicy Reels for Free - Exciting Bonus Features and Simple Gameplay" "" 0

# 2-5. "What we like" bullet list items (each preceded by an empty <w:r/>)
Replace-TextPreservingRuns "Intuitive user interface and controls" "Fantastic bonus features" "" 0
Replace-TextPreservingRuns "Suitable for both high rollers and beginners" "Simple rules and user interface" "" 0
Replace-TextPreservingRuns "Interesting bonus features" "Wide appeal to different player levels" "" 0
Replace-TextPreservingRuns "Appealing retro icons" "Variety of fruit symbols and fixed paylines" "" 0

# 6-7. "What we don't like" bullet list items (each preceded by an empty <w:r/>)
Replace-TextPreservingRuns "Limited selection of Symbols" "Limited number of paylines (20)" "" 0
Replace-TextPreservingRuns "Only 20 fixed paylines" "No progressive jackpot feature" "" 0

# 8. Bold meta-title run near the end (preceded by an empty <w:r/>), preserve bold formatting.
# Search starting after position 100 to skip the first (Heading1) occurrence.
Replace-TextPreservingRuns "Play Juicy Reels for Free - Exciting Bonus Features" "Play Juicy Reels for Free - Exciting Bonus Features and Simple Gameplay" "<w:rPr><w:b/></w:rPr>" 100

# 9. Italic meta-description run (preceded by an empty <w:r/>), preserve italic formatting
Replace-TextPreservingRuns "Discover Juicy Reels - a 6-reel 20-payline slot game with Wilds, Scatters and Free Spins. Suitable for all levels. Play free and experience unique features." "Experience the thrill of Juicy Reels with its fantastic bonus features and simple gameplay. Play for free now!" "<w:rPr><w:i/></w:rPr>" 0
